$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows for "Huesca" and "Huelva" were reordered (Huelva now listed first),
# taking their "Casos activos" (column C) figures along with their city names.
$ws.Range("A53").Value = "Huelva"
$ws.Range("C53").Value = 72

$ws.Range("A54").Value = "Huesca"
$ws.Range("C54").Value = 0

# Refresh the "Datos actualizados" timestamp shown in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 08:16"
